# Encoder BOM update
#  - fills in the parts list (switch, encoder, LED ring) with qty / price / links
#  - adds a "Link" hyperlink for the pre-existing MCU row
#  - restyles the header row (bigger bold font, centered)
#  - widens the Link column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header row styling ---------------------------------------------------
# Build the new header look on A1 first (its border matches B1:D1), then fan
# it out with a format-only paste; E1 gets its own pass since it sits on a
# different (inner) border variant.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 16
$ws.Range("A1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1").VerticalAlignment = -4108     # xlCenter

$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").Font.Size = 16
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").VerticalAlignment = -4108

$excel.CutCopyMode = $false
$ws.Rows(1).RowHeight = 29.4

# ---- widen the Link column -------------------------------------------------
$ws.Columns(5).ColumnWidth = 43.6

# ---- give the whole Link column the same banded fill/border as column A ---
# (covers the blank trailing rows too; the hyperlinks added below will only
# overlay the font, keeping this fill/border underneath)
$ws.Range("A2:A20").Copy() | Out-Null
$ws.Range("E2:E20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- row 3: Cherry MX switch -----------------------------------------------
$ws.Cells.Item(3, 1).Value = "CHERRY MX 540-MX1A-E1NW"
$ws.Cells.Item(3, 2).Value = "Mechanical pushbutton switch"
$ws.Cells.Item(3, 3).Value = 9
$ws.Cells.Item(3, 4).Value = 0.958

# ---- row 4: Wurth incremental encoder --------------------------------------
$ws.Cells.Item(4, 1).Value = "WURTH ELEKTRONIK 482009514001"
$ws.Cells.Item(4, 2).Value = "Mechanical incremental encoder with switch"
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 4.36

# ---- links for the switch / encoder rows -----------------------------------
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.mouser.com/ProductDetail/Cherry/MX1A-E1NW", "", "", "cherry mx - mouser.com")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.mouser.com/ProductDetail/Wurth-Elektronik/482009514001", "", "", "wurth incremental encoder - mouser.com")

# ---- row 5: WS2812B RGB LED ring -------------------------------------------
$ws.Cells.Item(5, 1).Value = "LED RGB WS2812B RING"
$ws.Cells.Item(5, 2).Value = "Ring of 8 RGB LED diodes"
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 2.5

$ws.Hyperlinks.Add($ws.Range("E5"), "https://botland.pl/pl/diody-led-rgb-i-rgbw/2943-pierscien-ws2812-8-x-diody-rgb-led.html", "", "", "ws2812 RGB LED - botland.pl")

# ---- link for the already-present MCU row (row 2) --------------------------
$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.lcsc.com/product-detail/STM32F401RCT6.html", "", "", "STM32F401RCT - lcsc.com")

# ---- selection, matching the saved workbook view ---------------------------
$ws.Range("E4").Select() | Out-Null
